# Update the dSF column (F) values for a handful of rows, per the
# "repull data, push all data, mean calculation" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -2
$ws.Range("F8").Value = -3
$ws.Range("F9").Value = 1
